$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where only the Price (D) and/or Volume(1h) (E) values change.
# D=$null means the Price column is unchanged for that row.
$updates = @(
    @{ Row=2; D="67.462.46"; E="  +0.62%  " },
    @{ Row=3; D="2.634.82"; E="  +0.83%  " },
    @{ Row=4; D=$null; E="  -0.04%  " },
    @{ Row=5; D="602.18"; E="  +1.46%  " },
    @{ Row=6; D="154.19"; E="  +0.32%  " },
    @{ Row=7; D=$null; E="  -0.04%  " },
    @{ Row=8; D="0.551"; E="  +1.64%  " },
    @{ Row=9; D="2.633.12"; E="  +0.80%  " },
    @{ Row=10; D="0.124"; E="  +6.64%  " },
    @{ Row=11; D=$null; E="  +0.64%  " },
    @{ Row=12; D=$null; E="  +0.22%  " },
    @{ Row=13; D="0.352"; E="  -0.69%  " },
    @{ Row=14; D="28.02"; E="  +1.51%  " },
    @{ Row=15; D="3.118.40"; E="  +0.85%  " },
    @{ Row=16; D=$null; E="  +1.42%  " },
    @{ Row=17; D="67.470.38"; E="  +0.71%  " },
    @{ Row=18; D="2.636.02"; E="  +0.77%  " },
    @{ Row=19; D="11.27"; E="  +0.16%  " },
    @{ Row=20; D="364.76"; E="  +1.51%  " },
    @{ Row=21; D="7.62"; E="  -3.69%  " },
    @{ Row=22; D=$null; E="  -0.39%  " },
    @{ Row=23; D=$null; E="  +6.40%  " },
    @{ Row=24; D="0.999"; E="  -0.10%  " },
    @{ Row=25; D="10.06"; E="  -1.77%  " },
    @{ Row=26; D="66.16"; E="  -7.80%  " },
    @{ Row=29; D="581.32"; E="  -6.79%  " },
    @{ Row=30; D=$null; E="  -3.02%  " },
    @{ Row=31; D="1.42"; E="  -2.29%  " },
    @{ Row=32; D="7.91"; E="  -1.03%  " },
    @{ Row=33; D=$null; E="  -0.10%  " },
    @{ Row=34; D="0.130"; E="  -1.74%  " },
    @{ Row=35; D="0.999"; E="  +0.00%  " },
    @{ Row=36; D="1.54"; E="  -1.34%  " },
    @{ Row=37; D="4.96"; E="  -0.45%  " },
    @{ Row=38; D="158.44"; E="  +3.05%  " },
    @{ Row=39; D="19.48"; E="  +0.45%  " },
    @{ Row=40; D="0.371"; E="  +0.56%  " },
    @{ Row=41; D="5.31"; E="  -3.30%  " },
    @{ Row=42; D="1.83"; E="  +0.80%  " },
    @{ Row=43; D="2.63"; E="  +1.78%  " },
    @{ Row=44; D="41.23"; E="  -0.25%  " },
    @{ Row=45; D=$null; E="  +0.00%  " },
    @{ Row=46; D="16.36"; E="  -0.74%  " },
    @{ Row=47; D="156.36"; E="  +0.65%  " },
    @{ Row=48; D="0.0₆0289"; E="  -3.43%  " },
    @{ Row=49; D=$null; E="  -0.87%  " }
)

foreach ($u in $updates) {
    if ($null -ne $u.D) {
        # Force the Price cell to a text value so strings like "602.18" or
        # "67.462.46" are not silently coerced into floating point numbers
        # (which would lose precision / change representation) or scientific
        # notation (e.g. "0.0000104").
        $ws.Cells.Item($u.Row, 4).NumberFormat = "@"
        $ws.Cells.Item($u.Row, 4).Value = $u.D
        $ws.Cells.Item($u.Row, 4).Style = "Normal"
    }
    $ws.Cells.Item($u.Row, 5).Value = $u.E
}

# Rows 27 and 28 swap order (WrappedeETH now ranks above PEPE) and get new
# Price / Volume(1h) figures.
$ws.Cells.Item(27, 2).Value = "WrappedeETH"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.769.82"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +0.71%  "

$ws.Cells.Item(28, 2).Value = "PEPE"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.0000104"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +1.11%  "

# Rows 50 and 51 swap order (Mantle now ranks above InjectiveProtocol) and
# get new Price / Volume(1h) figures.
$ws.Cells.Item(50, 2).Value = "Mantle"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "0.628"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  +0.62%  "

$ws.Cells.Item(51, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "20.97"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  +0.15%  "
